$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 5937
$ws.Range("I3").Value = 6196
$ws.Range("I4").Value = 1419
$ws.Range("I5").Value = 575
$ws.Range("I6").Value = 7003
$ws.Range("I7").Value = 21130

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 164
$ws.Range("I4").Value = 85
$ws.Range("I6").Value = 151
$ws.Range("I8").Value = 1267
$ws.Range("I9").Value = 101
$ws.Range("I10").Value = 149
$ws.Range("I11").Value = 318
$ws.Range("I15").Value = 244
$ws.Range("I19").Value = 588
$ws.Range("I20").Value = 527
$ws.Range("I24").Value = 60
$ws.Range("I25").Value = 113
$ws.Range("I27").Value = 189
$ws.Range("I29").Value = 1310
$ws.Range("I31").Value = 210
$ws.Range("I33").Value = 956
$ws.Range("I36").Value = 284
$ws.Range("I37").Value = 677
$ws.Range("I41").Value = 92
$ws.Range("I42").Value = 724
$ws.Range("I44").Value = 158
$ws.Range("I48").Value = 285
$ws.Range("I49").Value = 143
$ws.Range("I50").Value = 104
$ws.Range("I52").Value = 458
$ws.Range("I53").Value = 222
$ws.Range("I54").Value = 433
$ws.Range("I55").Value = 231
$ws.Range("I63").Value = 68
$ws.Range("I65").Value = 493
$ws.Range("I66").Value = 62
$ws.Range("I67").Value = 817
$ws.Range("I69").Value = 47
$ws.Range("I73").Value = 193
$ws.Range("I75").Value = 68
$ws.Range("I76").Value = 303
$ws.Range("I79").Value = 600
$ws.Range("I83").Value = 455
$ws.Range("I85").Value = 962
$ws.Range("I87").Value = 50
$ws.Range("I88").Value = 192
$ws.Range("I90").Value = 256
$ws.Range("I92").Value = 59
$ws.Range("I93").Value = 120
$ws.Range("I94").Value = 221
$ws.Range("I95").Value = 323
$ws.Range("I97").Value = 177
$ws.Range("I99").Value = 379
$ws.Range("I101").Value = 21130

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 272
$ws.Range("I3").Value = 370
$ws.Range("I6").Value = 241
$ws.Range("I7").Value = 962

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("I6").Value = 16
$ws.Range("I7").Value = 47

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 162
$ws.Range("I4").Value = 38
$ws.Range("I6").Value = 119
$ws.Range("I7").Value = 458

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("I2").Value = 132
$ws.Range("I3").Value = 60
$ws.Range("I6").Value = 84
$ws.Range("I7").Value = 318

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 387
$ws.Range("I3").Value = 360
$ws.Range("I5").Value = 37
$ws.Range("I7").Value = 1267

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I3").Value = 46
$ws.Range("I6").Value = 102
$ws.Range("I7").Value = 222

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 203
$ws.Range("I3").Value = 226
$ws.Range("I6").Value = 195
$ws.Range("I7").Value = 677

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 139
$ws.Range("I7").Value = 379

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 190
$ws.Range("I3").Value = 300
$ws.Range("I6").Value = 251
$ws.Range("I7").Value = 817

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I2").Value = 65
$ws.Range("I7").Value = 210

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I2").Value = 163
$ws.Range("I3").Value = 151
$ws.Range("I7").Value = 493

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 156
$ws.Range("I7").Value = 455

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 111
$ws.Range("I3").Value = 117
$ws.Range("I7").Value = 323

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I2").Value = 211
$ws.Range("I3").Value = 363
$ws.Range("I6").Value = 302
$ws.Range("I7").Value = 956

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I2").Value = 27
$ws.Range("I7").Value = 143

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 207
$ws.Range("I7").Value = 433

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I3").Value = 452
$ws.Range("I6").Value = 363
$ws.Range("I7").Value = 1310

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 198
$ws.Range("I3").Value = 180
$ws.Range("I6").Value = 175
$ws.Range("I7").Value = 588

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("I2").Value = 49
$ws.Range("I3").Value = 47
$ws.Range("I4").Value = 14
$ws.Range("I5").Value = 3
$ws.Range("I7").Value = 158

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I4").Value = 34
$ws.Range("I7").Value = 285

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I2").Value = 58
$ws.Range("I3").Value = 69
$ws.Range("I6").Value = 140
$ws.Range("I7").Value = 303

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("I5").Value = 4
$ws.Range("I7").Value = 151

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I3").Value = 30
$ws.Range("I6").Value = 22
$ws.Range("I7").Value = 92

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 233
$ws.Range("I6").Value = 236
$ws.Range("I7").Value = 724

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I6").Value = 65
$ws.Range("I7").Value = 149

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 231

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("I3").Value = 24
$ws.Range("I7").Value = 60

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I2").Value = 175
$ws.Range("I3").Value = 194
$ws.Range("I7").Value = 600

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I3").Value = 152
$ws.Range("I6").Value = 183
$ws.Range("I7").Value = 527

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I2").Value = 83
$ws.Range("I3").Value = 92
$ws.Range("I5").Value = 10
$ws.Range("I6").Value = 89
$ws.Range("I7").Value = 284

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I6").Value = 50
$ws.Range("I7").Value = 120

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("I6").Value = 128
$ws.Range("I7").Value = 221

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I2").Value = 43
$ws.Range("I3").Value = 31
$ws.Range("I7").Value = 113

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 55
$ws.Range("I5").Value = 11
$ws.Range("I6").Value = 90
$ws.Range("I7").Value = 244

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I2").Value = 29
$ws.Range("I7").Value = 104

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("I2").Value = 20
$ws.Range("I4").Value = 7
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("I3").Value = 34
$ws.Range("I7").Value = 101

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I2").Value = 63
$ws.Range("I3").Value = 60
$ws.Range("I7").Value = 193

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I2").Value = 57
$ws.Range("I7").Value = 164

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I6").Value = 114
$ws.Range("I7").Value = 177

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("I3").Value = 11
$ws.Range("I6").Value = 24
$ws.Range("I7").Value = 59

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("I2").Value = 56
$ws.Range("I7").Value = 192

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I3").Value = 40
$ws.Range("I7").Value = 189

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("I6").Value = 17
$ws.Range("I7").Value = 68

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 84
$ws.Range("I7").Value = 256

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I6").Value = 28
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 50
